$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 4) with a new admin account "putong" in column A,
# following the existing "qishilong" (A2) / "admin" (A3) entries.
$ws.Range("A4").Value = "putong"
